# faturamento_diario.xlsx update:
# A new daily-sales record (day 10 of July/2025) was added to the dataset.
# It belongs right after the existing July rows (day 1..9, sheet rows 2-10)
# and before the first June row (previously row 11), so insert a new row at
# row 11 - pushing the June/May/April rows down by one - and populate it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 11 (shifts old rows 11..101 down to 12..102,
# same as right-clicking row 11's header and choosing "Insert" in Excel).
$ws.Rows.Item(11).Insert()

# Fill in the new row with the July (month 7 / 2025) record for day 10.
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 21266
$ws.Range("C11").Value = 7
$ws.Range("D11").Value = 2025
$ws.Range("E11").Value = "07/2025"
